# Append a new "(2025, 2035]" decade-bucket row (row 11) to the page-counts
# table, mirroring the all-zero/empty stats produced for a bucket with no
# matching pages (per the BERTopic re-summarization run).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row values.
$ws.Cells.Item(11, 1).Value = "(2025, 2035]"
$ws.Cells.Item(11, 2).Value = 0

# Match the existing table formatting: A11 picks up the bold/border/
# center-top style used by the other "decade" cells (A2:A10), and C11/D11
# pick up the plain numeric-column formatting (no explicit style) used by
# the rest of the mean/std columns -- even though they carry no value for
# this bucket, so they still show up as real (empty) cells.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)

$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial(-4122)

$ws.Range("D10").Copy()
$ws.Range("D11").PasteSpecial(-4122)

Write-Output "Added row 11 to $($ws.Name); used range now $($ws.UsedRange.Address())"
